$wb = $excel.ActiveWorkbook

# --- "Rushing" sheet ---
$wsRushing = $wb.Worksheets.Item("Rushing")

# Update P.Barber's rushing stats (1DATT/2DATT) before removing K.Drake's row
$wsRushing.Cells.Item(7, 3).Value = 30
$wsRushing.Cells.Item(7, 4).Value = 12

# Remove K.Drake from the "Rushing" sheet (row 6)
$wsRushing.Rows.Item(6).Delete()

# Fix up the sequential index values in column A (rows 2..12) after the row shift
for ($i = 2; $i -le 12; $i++) {
    $wsRushing.Cells.Item($i, 1).Value = $i - 2
}

# --- "Receiving" sheet ---
$wsReceiving = $wb.Worksheets.Item("Receiving")

# Update P.Barber's receiving stats (Short Target/Short Comp) before removing K.Drake's row
$wsReceiving.Cells.Item(4, 3).Value = 10
$wsReceiving.Cells.Item(4, 4).Value = 8

# Remove K.Drake from the "Receiving" sheet (row 3)
$wsReceiving.Rows.Item(3).Delete()

# Fix up the sequential index values in column A (rows 2..13) after the row shift
for ($i = 2; $i -le 13; $i++) {
    $wsReceiving.Cells.Item($i, 1).Value = $i - 2
}
